$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.874.72"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.855.21"
$ws.Range("E3").Value = "  +3.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.20"
$ws.Range("E6").Value = "  -2.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.851.67"
$ws.Range("E7").Value = "  +2.99%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("E11").Value = "  -2.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("E13").Value = "  -3.07%  "
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.497.04"
$ws.Range("E15").Value = "  +2.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.857.14"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.039.62"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.54"
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.49"
$ws.Range("E19").Value = "  +4.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.17"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "485.09"
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("E24").Value = "  +4.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.75"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.09"
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.96"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.93"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.005.53"
$ws.Range("E32").Value = "  +3.08%  "
$ws.Range("E33").Value = "  -3.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.24"
$ws.Range("E34").Value = "  +2.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.802.88"
$ws.Range("E35").Value = "  +3.38%  "
$ws.Range("E36").Value = "  -1.77%  "
$ws.Range("E37").Value = "  +1.67%  "
$ws.Range("E38").Value = "  +3.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.89"
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "439.67"
$ws.Range("E42").Value = "  +1.89%  "
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.48"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.40"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.70"
$ws.Range("E48").Value = "  +12.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.73"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.833.13"
$ws.Range("E50").Value = "  +1.76%  "
$ws.Range("E51").Value = "  +1.72%  "
